# Insert a new weekly record as row 17, pushing the existing rows 17-55
# down to 18-56 (dimension grows from A1:R55 to A1:R56).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 17..55 down by one to make room for the new record.
$ws.Rows.Item(17).Insert()

# New record for row 17.
$ws.Cells.Item(17, 1).Value = 10
$ws.Cells.Item(17, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(17, 3).Value = "La Araucanía"
$ws.Cells.Item(17, 4).Value = 44575
$ws.Cells.Item(17, 5).Value = 9
$ws.Cells.Item(17, 6).Value = 100112030
$ws.Cells.Item(17, 7).Value = "Poroto granado"
$ws.Cells.Item(17, 8).Value = "Sin especificar"
$ws.Cells.Item(17, 9).Value = "Primera"
$ws.Cells.Item(17, 10).Value = 125
$ws.Cells.Item(17, 11).Value = 28000
$ws.Cells.Item(17, 12).Value = 28000
$ws.Cells.Item(17, 13).Value = 28000
$ws.Cells.Item(17, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(17, 15).Value = "Región del Maule"
$ws.Cells.Item(17, 16).Value = 1120
$ws.Cells.Item(17, 17).Value = 25
$ws.Cells.Item(17, 18).Value = "Hortaliza"
